# Apply the corrected panel-analysis results to the table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Prevalence" row label (drop the trailing units note)
$ws.Range("A2").Value = "Prevalence"

# Updated coefficient (col4) values
$ws.Range("B3").Value = -0.00049
$ws.Range("B4").Value = -0.00304
$ws.Range("B5").Value = -0.0158
$ws.Range("B6").Value = 0.01254

# Updated 95% confidence intervals (ci95)
$ws.Range("C3").Value = "(-0.00101 - 2e-05)"
$ws.Range("C4").Value = "(-0.00424 - -0.00185)"
$ws.Range("C5").Value = "(-0.01846 - -0.01313)"
$ws.Range("C6").Value = "(0.00844 - 0.01664)"
$ws.Range("C7").Value = "(0.00224 - 0.00293)"
$ws.Range("C8").Value = "(-0.09183 - -0.07677)"
